$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'290.95"
$ws.Range("E2").Value = "'-3.97%"
$ws.Range("D3").Value = "'30.97"
$ws.Range("E3").Value = "'-3.67%"
$ws.Range("D4").Value = "'4.945"
$ws.Range("E4").Value = "'0.40%"
$ws.Range("D5").Value = "'0.07162"
$ws.Range("E5").Value = "'-8.57%"
$ws.Range("D6").Value = "'1.794"
$ws.Range("E6").Value = "'-11.84%"
$ws.Range("E7").Value = "'-2.10%"
$ws.Range("D8").Value = "'3.734"
$ws.Range("E8").Value = "'-2.57%"
$ws.Range("E9").Value = "'-2.99%"
$ws.Range("D10").Value = "'0.1648"
$ws.Range("E10").Value = "'-6.47%"
$ws.Range("D11").Value = "'0.07697"
$ws.Range("E11").Value = "'-2.31%"
$ws.Range("D12").Value = "'0.08115"
$ws.Range("E12").Value = "'-5.53%"
$ws.Range("D13").Value = "'0.03019"
$ws.Range("E13").Value = "'-4.55%"
$ws.Range("D14").Value = "'0.1004"
$ws.Range("E14").Value = "'-0.08%"
$ws.Range("D15").Value = "'0.001491"
$ws.Range("E15").Value = "'-1.32%"
$ws.Range("D16").Value = "'0.005737"
$ws.Range("E16").Value = "'-2.77%"
$ws.Range("D17").Value = "'3.475"
$ws.Range("E17").Value = "'0.31%"
$ws.Range("D18").Value = "'2.081"
$ws.Range("E18").Value = "'-3.49%"
$ws.Range("D19").Value = "'0.3279"
$ws.Range("E19").Value = "'0.06%"
$ws.Range("E20").Value = "'-3.41%"
$ws.Range("D21").Value = "'4.039"
$ws.Range("E21").Value = "'-5.50%"
$ws.Range("D22").Value = "'0.1996"
$ws.Range("E22").Value = "'-0.14%"
$ws.Range("D23").Value = "'0.04516"
$ws.Range("E23").Value = "'-1.24%"
$ws.Range("D24").Value = "'0.001211"
$ws.Range("E24").Value = "'-1.01%"
$ws.Range("D25").Value = "'0.004010"
$ws.Range("E25").Value = "'-9.85%"
$ws.Range("E26").Value = "'-0.10%"
$ws.Range("D39").Value = "'0.01607"
$ws.Range("E39").Value = "'-7.76%"
$ws.Range("D40").Value = "'0.04386"
$ws.Range("E40").Value = "'-8.28%"
$ws.Range("D41").Value = "'0.007386"
$ws.Range("E41").Value = "'-2.29%"
$ws.Range("D42").Value = "'0.1307"
$ws.Range("E42").Value = "'-4.22%"
$ws.Range("D43").Value = "'0.001997"
$ws.Range("E43").Value = "'-15.41%"
$ws.Range("D44").Value = "'0.009222"
$ws.Range("E44").Value = "'-12.67%"
$ws.Range("D45").Value = "'0.00005923"
$ws.Range("E45").Value = "'-6.23%"
$ws.Range("E46").Value = "'-0.19%"
$ws.Range("D47").Value = "'2.246"
$ws.Range("E47").Value = "'172.73%"
$ws.Range("D48").Value = "'0.002995"
$ws.Range("E48").Value = "'-3.43%"
$ws.Range("E49").Value = "'-0.19%"
$ws.Range("E50").Value = "'-0.19%"
